$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - column F ("想去人数") updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 691
$ws1.Range("F5").Value = 579
$ws1.Range("F6").Value = 324
$ws1.Range("F7").Value = 2833
$ws1.Range("F9").Value = 8077
$ws1.Range("F10").Value = 207
$ws1.Range("F11").Value = 477
$ws1.Range("F13").Value = 409
$ws1.Range("F14").Value = 49

# Sheet "全部类型" (sheet4) - column F ("想去人数") updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 691
$ws4.Range("F5").Value = 579
$ws4.Range("F6").Value = 324
$ws4.Range("F9").Value = 2833
$ws4.Range("F11").Value = 8077
$ws4.Range("F12").Value = 207
$ws4.Range("F13").Value = 477
$ws4.Range("F17").Value = 409
$ws4.Range("F18").Value = 49
